$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (工兵机器人 / Miner bot) - effect text (column E): add "可以" (optional) before
# "将主牌堆第1张陷阱牌加入手牌"
$ws.Cells.Item(3, 5).Value = "有牌进入弃牌堆时：点数减1，然后如果本牌①与玩家敌对，则将主牌堆第1张陷阱牌放在房间区任意位置。②受玩家控制，则可以将主牌堆第1张陷阱牌加入手牌，或将房间区1张陷阱牌送墓。"

# Row 4 (迷你工厂 / Mini factory) - effect text (column E): rewrite the robot-retrieval effect
$ws.Cells.Item(4, 5).Value = "有牌进入弃牌堆时：点数减1，然后如果本牌①在房间区，则将墓地第1张“机器人”牌放在本牌前方相邻的单元。②在手牌，则可以将墓地第1张“机器人”牌加入手牌。"

# Row 5 (自爆机器人 / Kamikaze bot) - effect text (column E): add "可以" (optional)
$ws.Cells.Item(5, 5).Value = "有牌进入弃牌堆时：点数减1。<br>`n点数为0时：如果本牌①与玩家敌对，则玩家受到1伤害。②受玩家控制，则可以选场上1张牌送墓。"

# Row 6 (哨戒机器人 / Sentinel bot) - effect text (column E): add "可以" (optional)
$ws.Cells.Item(6, 5).Value = "有牌进入弃牌堆时：点数减1。<br>`n回合结束时：如果本牌①与玩家敌对，则将主牌堆第1张机器人牌放在房间区任意位置。②受玩家控制，则可以将主牌堆第1张机器人牌加入手牌。"

# Row 15 (黏菌 / Goo) - effect text (column E): completely rewritten effect
$ws.Cells.Item(15, 5).Value = "送墓时：选手牌或房间区1张“史莱姆”牌，其点数加1。"

# Reflect the cursor/viewport position recorded in the saved file
$ws.Rows.Item(15).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(6).AutoFit()

$ws.Application.Goto($ws.Range("E7"), $true)
$ws.Range("E7").Select()
